$d = $word.ActiveDocument

# New entries for the 12/9/10 meeting, appended after the last paragraph
# ("12/6/10: Cleaned names from ish dataset") and before the section break.
$entries = @(
    @{ Style = "Heading2"; Text = "12/9/10:" },
    @{ Style = "Heading3"; Text = "Progress to point – fixed typedefs, fixed includes stack to start with initial file, put comments in their own start states but need to fix line numbers" },
    @{ Style = "Heading3"; Text = "Meeting:" },
    @{ Style = "Heading4"; Text = "Includes:" },
    @{ Style = "Heading5"; Text = "All std stuff is in /usr/include" },
    @{ Style = "Heading5"; Text = "If file no found in local folders search/try in /usr/include (which matches the definition for inclusion on “” files, <> is just a suppression of the local folder search)" },
    @{ Style = "Heading5"; Text = "Try to run on stdlib (io, strings) and if epic fails then see the cost of fixing it. If not worth fixing keep hack of special size_t in lexer" },
    @{ Style = "Heading4"; Text = "atexit(funcptr) calls function at program exit" },
    @{ Style = "Heading4"; Text = "linenumbers and filenames not critical" }
)

foreach ($entry in $entries) {
    $count = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs($count)
    $rng = $lastPara.Range
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()

    $newCount = $d.Paragraphs.Count
    $newPara = $d.Paragraphs($newCount)
    $newPara.Range.Text = $entry.Text
    $newPara.Style = $entry.Style
}
